# Tidsregistrering - Dan Toft: log a few more days of hours (GUI work,
# "alt andet", bugfixing) on the "Tidsregistrering" sheet's time table.
# The sheet has running formulas for Dato/Starttid/Timer, so we only need
# to fill in the Aktivitet (D) and Sluttid (F) cells for rows 34-37; the
# dependent formula columns (B, E, G) recompute on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")
$ws.Activate()

# Row 34: finished the "GUI" activity at 15:30 (F >= 15:30 rolls the next
# row's date forward by one day via the existing B-column formula).
$ws.Range("D34").Value = "GUI"
$ws.Range("F34").Value = 0.64583333333333337

# Row 35: "Alt andet end gui og alligevel endte jeg med at lave GUI",
# also ending 15:30.
$ws.Range("D35").Value = "Alt andet end gui og alligevel endte jeg med at lave GUI"
$ws.Range("F35").Value = 0.64583333333333337

# Row 36: same activity text/formatting as row 35 - copy the cell so the
# (slightly different) banding style on D35 carries over to D36 too.
$ws.Range("D35").Copy($ws.Range("D36"))
$ws.Range("F36").Value = 0.64583333333333337

# Row 37: "Bugfixing", ending earlier at 14:00 (no date rollover for row 38).
$ws.Range("D37").Value = "Bugfixing"
$ws.Range("F37").Value = 0.58333333333333337

# Restore the selection/scroll position left behind by the edit session.
$ws.Range("F38").Select()
